$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 5 format (the "separator" border-group style 6/7) onto row 12
$ws.Range("A5:E5").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A12").Value = $null

# Copy row 2 format (style 4/5 group body, row height 43.2 non-custom) onto row 13
$ws.Range("A2:E2").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)  # xlPasteFormats

# Set values in the shared-string-append order matching the target workbook:
# C13 first (-> index 42), A13 next (-> index 43), D13 (-> 44), E13 (-> 45)
$ws.Range("C13").Value = " Team [team:] is in charge!"
$ws.Range("A13").Value = "SCRIPT/G01P03A/um2406.ssb"
$ws.Range("D13").Value = " Вперёд, Команда [team:]!"
$ws.Range("E13").Value = " Âðåñæä, Ëïíàîäà [team:]!"

$ws.Range("B13").Clear()
$ws.Rows.Item(13).RowHeight = 43.2

$ws.Range("D13").Select() | Out-Null
